$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 447 - this shifts existing rows 447:508
# down to 448:509 (matching the dimension change A1:R508 -> A1:R509).
$ws.Rows.Item(447).Insert()

# Populate the newly-inserted row 447 with its data.
$ws.Range("A447").Value = 3
$ws.Range("B447").Value = 'Femacal de La Calera'
$ws.Range("C447").Value = 'Coquimbo'
$ws.Range("D447").Value = 44984
$ws.Range("E447").Value = 5
$ws.Range("F447").Value = 100112043
$ws.Range("G447").Value = 'Pepino ensalada'
$ws.Range("H447").Value = 'Sin especificar'
$ws.Range("I447").Value = 'Primera'
$ws.Range("J447").Value = 135
$ws.Range("K447").Value = 6000
$ws.Range("L447").Value = 6500
$ws.Range("M447").Value = 6259
$ws.Range("N447").Value = '$/caja 60 unidades'
$ws.Range("O447").Value = 'Región de Arica y Parinacota'
$ws.Range("P447").Value = 104
$ws.Range("Q447").Value = 60
$ws.Range("R447").Value = 'Hortaliza'
